$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle_D2 = $ws.Range('D2').Style
$ws.Range('D2').Value = "'44.688.33"
$ws.Range('D2').Style = $origStyle_D2
$ws.Range('E2').Value = '  +4.06%  '
$origStyle_D3 = $ws.Range('D3').Style
$ws.Range('D3').Value = "'2.424.08"
$ws.Range('D3').Style = $origStyle_D3
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('E4').Value = '  -0.06%  '
$origStyle_D5 = $ws.Range('D5').Style
$ws.Range('D5').Value = "'315.43"
$ws.Range('D5').Style = $origStyle_D5
$ws.Range('E5').Value = '  +4.01%  '
$origStyle_D6 = $ws.Range('D6').Style
$ws.Range('D6').Value = "'101.70"
$ws.Range('D6').Style = $origStyle_D6
$ws.Range('E6').Value = '  +6.67%  '
$ws.Range('E7').Value = '  +2.60%  '
$origStyle_D9 = $ws.Range('D9').Style
$ws.Range('D9').Value = "'0.524"
$ws.Range('D9').Style = $origStyle_D9
$ws.Range('E9').Value = '  +8.46%  '
$ws.Range('E10').Value = '  +3.99%  '
$origStyle_D11 = $ws.Range('D11').Style
$ws.Range('D11').Value = "'0.0801"
$ws.Range('D11').Style = $origStyle_D11
$origStyle_D12 = $ws.Range('D12').Style
$ws.Range('D12').Value = "'18.98"
$ws.Range('D12').Style = $origStyle_D12
$ws.Range('E12').Value = '  +2.93%  '
$ws.Range('E13').Value = '  -2.30%  '
$origStyle_D14 = $ws.Range('D14').Style
$ws.Range('D14').Value = "'6.97"
$ws.Range('D14').Style = $origStyle_D14
$ws.Range('E14').Value = '  +3.53%  '
$origStyle_D15 = $ws.Range('D15').Style
$ws.Range('D15').Value = "'2.802.83"
$ws.Range('D15').Style = $origStyle_D15
$ws.Range('E15').Value = '  +2.73%  '
$origStyle_D16 = $ws.Range('D16').Style
$ws.Range('D16').Value = "'2.403.39"
$ws.Range('D16').Style = $origStyle_D16
$ws.Range('E16').Value = '  +2.54%  '
$origStyle_D17 = $ws.Range('D17').Style
$ws.Range('D17').Value = "'0.834"
$ws.Range('D17').Style = $origStyle_D17
$ws.Range('E17').Value = '  +5.28%  '
$origStyle_D18 = $ws.Range('D18').Style
$ws.Range('D18').Value = "'44.558.67"
$ws.Range('D18').Style = $origStyle_D18
$ws.Range('E18').Value = '  +3.79%  '
$origStyle_D19 = $ws.Range('D19').Style
$ws.Range('D19').Value = "'12.44"
$ws.Range('D19').Style = $origStyle_D19
$ws.Range('E19').Value = '  +4.64%  '
$origStyle_D20 = $ws.Range('D20').Style
$ws.Range('D20').Value = "'6.41"
$ws.Range('D20').Style = $origStyle_D20
$ws.Range('E20').Value = '  +2.45%  '
$ws.Range('E21').Value = '  +4.38%  '
$origStyle_D22 = $ws.Range('D22').Style
$ws.Range('D22').Value = "'68.81"
$ws.Range('D22').Style = $origStyle_D22
$ws.Range('E22').Value = '  +1.09%  '
$origStyle_D23 = $ws.Range('D23').Style
$ws.Range('D23').Value = "'242.84"
$ws.Range('D23').Style = $origStyle_D23
$ws.Range('E23').Value = '  +3.36%  '
$ws.Range('E24').Value = '  +5.94%  '
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('E28').Value = '  -3.58%  '
$origStyle_D29 = $ws.Range('D29').Style
$ws.Range('D29').Value = "'9.52"
$ws.Range('D29').Style = $origStyle_D29
$ws.Range('E29').Value = '  +1.93%  '
$origStyle_D30 = $ws.Range('D30').Style
$ws.Range('D30').Value = "'33.50"
$ws.Range('D30').Style = $origStyle_D30
$ws.Range('E30').Value = '  +4.03%  '
$origStyle_D31 = $ws.Range('D31').Style
$ws.Range('D31').Value = "'48.27"
$ws.Range('D31').Style = $origStyle_D31
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('E32').Value = '  +21.19%  '
$origStyle_D33 = $ws.Range('D33').Style
$ws.Range('D33').Value = "'19.39"
$ws.Range('D33').Style = $origStyle_D33
$ws.Range('E33').Value = '  +10.73%  '
$origStyle_D34 = $ws.Range('D34').Style
$ws.Range('D34').Value = "'0.0781"
$ws.Range('D34').Style = $origStyle_D34
$ws.Range('E34').Value = '  +9.29%  '
$origStyle_D35 = $ws.Range('D35').Style
$ws.Range('D35').Value = "'5.18"
$ws.Range('D35').Style = $origStyle_D35
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  +2.91%  '
$origStyle_D38 = $ws.Range('D38').Style
$ws.Range('D38').Value = "'4.49"
$ws.Range('D38').Style = $origStyle_D38
$ws.Range('E38').Value = '  +4.23%  '
$ws.Range('E39').Value = '  +1.07%  '
$origStyle_D40 = $ws.Range('D40').Style
$ws.Range('D40').Value = "'121.00"
$ws.Range('D40').Style = $origStyle_D40
$ws.Range('E40').Value = '  -4.42%  '
$ws.Range('E41').Value = '  +2.02%  '
$ws.Range('E42').Value = '  -2.27%  '
$origStyle_D43 = $ws.Range('D43').Style
$ws.Range('D43').Value = "'21.04"
$ws.Range('D43').Style = $origStyle_D43
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('E44').Value = '  +4.76%  '
$origStyle_D45 = $ws.Range('D45').Style
$ws.Range('D45').Value = "'1.943.50"
$ws.Range('D45').Style = $origStyle_D45
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('E47').Value = '  +9.09%  '
$origStyle_D48 = $ws.Range('D48').Style
$ws.Range('D48').Value = "'9.49"
$ws.Range('D48').Style = $origStyle_D48
$ws.Range('E48').Value = '  +2.43%  '
$origStyle_D49 = $ws.Range('D49').Style
$ws.Range('D49').Value = "'1.69"
$ws.Range('D49').Style = $origStyle_D49
$ws.Range('E49').Value = '  +11.97%  '
$origStyle_D50 = $ws.Range('D50').Style
$ws.Range('D50').Value = "'54.83"
$ws.Range('D50').Style = $origStyle_D50
$ws.Range('E50').Value = '  +7.00%  '
$origStyle_D51 = $ws.Range('D51').Style
$ws.Range('D51').Value = "'75.68"
$ws.Range('D51').Style = $origStyle_D51
$ws.Range('E51').Value = '  +5.88%  '
